$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E3").Value = 16.696
$ws.Range("A4").Value = -20.63
$ws.Range("A6").Value = -22.255
$ws.Range("A7").Value = -19.975
$ws.Range("C7").Value = -12.489
$ws.Range("A8").Value = -22.114
$ws.Range("C11").Value = -12.924
$ws.Range("C12").Value = -11.926
$ws.Range("D12").Value = -7.832000000000001
$ws.Range("E12").Value = 17.344
$ws.Range("D13").Value = -8.166999999999998
$ws.Range("E13").Value = 16.473
$ws.Range("D14").Value = -7.817
$ws.Range("C15").Value = -13.232
$ws.Range("A16").Value = -21.82
$ws.Range("D16").Value = -8.666
$ws.Range("D19").Value = -8.017999999999999
$ws.Range("A20").Value = -19.884
$ws.Range("C20").Value = -12.203
$ws.Range("D20").Value = -7.386
$ws.Range("A21").Value = -20.048
$ws.Range("C21").Value = -12.293
$ws.Range("C22").Value = -13.137
$ws.Range("D22").Value = -7.887
$ws.Range("E22").Value = 16.871
$ws.Range("C23").Value = -12.223
$ws.Range("E25").Value = 17.179
$ws.Range("A28").Value = -21.914
$ws.Range("A29").Value = -21.305
$ws.Range("C29").Value = -12.696
$ws.Range("E29").Value = 17.073
$ws.Range("A30").Value = -21.954
$ws.Range("A32").Value = -21.681
$ws.Range("C34").Value = -11.818
$ws.Range("E34").Value = 16.994
$ws.Range("D36").Value = -7.922
$ws.Range("A40").Value = -20.132
$ws.Range("C42").Value = -12.287
$ws.Range("C43").Value = -13.24
$ws.Range("D43").Value = -7.874000000000001
$ws.Range("E43").Value = 16.696
$ws.Range("C44").Value = -12.74
$ws.Range("C45").Value = -13.054
$ws.Range("A46").Value = -21.859
$ws.Range("C46").Value = -12.664
$ws.Range("D46").Value = -8.336
$ws.Range("E48").Value = 17.022
$ws.Range("C50").Value = -12.609
$ws.Range("D50").Value = -8.431000000000001
$ws.Range("A51").Value = -21.649
$ws.Range("C51").Value = -11.177
$ws.Range("A52").Value = -21.988
$ws.Range("A57").Value = -21.875
$ws.Range("C57").Value = -13.033
$ws.Range("A59").Value = -22.057
$ws.Range("E60").Value = 16.386
$ws.Range("A62").Value = -22.103
$ws.Range("C65").Value = -12.4
$ws.Range("A66").Value = -21.668
$ws.Range("C66").Value = -11.047
$ws.Range("C67").Value = -11.417
$ws.Range("E68").Value = 17.337
$ws.Range("E70").Value = 17.46
$ws.Range("E71").Value = 17.036
$ws.Range("A73").Value = -20.509
$ws.Range("E73").Value = 16.654
$ws.Range("A74").Value = -21.112
$ws.Range("D76").Value = -7.523000000000001
$ws.Range("A77").Value = -20.371
$ws.Range("E78").Value = 16.669
$ws.Range("C79").Value = -12.663
$ws.Range("C84").Value = -13.424
$ws.Range("C87").Value = -13.101
$ws.Range("E87").Value = 16.319
$ws.Range("A92").Value = -21.766
$ws.Range("C92").Value = -11.434
$ws.Range("E92").Value = 17.517
$ws.Range("D95").Value = -7.536000000000001
$ws.Range("C97").Value = -12.038
$ws.Range("D97").Value = -8.161
$ws.Range("D99").Value = -7.711
$ws.Range("A100").Value = -22.122
$ws.Range("E101").Value = 16.54
